$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D32").Value = "Graph 유형 정리 (GNN)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/389"

$ws.Range("D36").Value = "Anomaly detection using imaging of time series data"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/380"

$ws.Range("D46").Value = "Labor Induction (유도분만)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/487"

$ws.Range("D50").Value = "하사비스, 결국 노벨상 수상으로 가나?"
$ws.Range("E50").Value = "http://incredible.egloos.com/7548675"
